$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new values (columns A..AH) ---
$row2 = @(45048.50694444445,9.173,6.534,2.853,20.581,14.857,6.617,20.037,11.575,4.738,6.004,8.497,8.644,2.821,7.519,9.989000000000001,7.075,2.196,0.759,107.045,20.658,6.941,12.856,7.392,1.396,11.913,6.131,5.786,6.628,8.785,2.215,18.053,3.621,8.676)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @(45048.51388888889,16.071,11.926,1.641,35.495,28.23,12.224,46.213,19.724,8.73,12.285,14.446,15.013,4.398,12.783,17.91,11.141,1.142,0.665,187.236,35.696,11.799,23.585,12.612,1.782,24.013,10.422,9.391,11.004,14.978,1.023,42.418,6.506,14.751)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @(45048.52083333334,24.313,18.179,1.523,53.289,43.324,18.807,73.226,29.617,13.274,19.222,21.544,22.593,6.389,19.174,27.174,16.275,0.891,0.8139999999999999,284.54,53.607,17.698,35.912,18.971,2.562,36.579,15.633,13.892,16.323,22.48,0.671,66.821,9.927,22.127)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @(45048.52777777778,6.58,4.88,0.74,14.63,11.55,4.91,26.22,8.1,3.68,4.94,6,6.19,1.89,5.26,7.4,4.66,0.58,0.25,72.83,14.98,4.86,9.869999999999999,5.18,0.72,12.59,4.29,3.91,4.57,6.17,0.5,24.42,2.64,6.08)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

# --- Delete row 6 (dataset now has one fewer row) ---
$ws.Rows.Item(6).Delete()

# --- Widen columns J, K, L, V, Z, AB from 7 to 8 characters ---
$wideCols = @(10, 11, 12, 22, 26, 28)
foreach ($c in $wideCols) {
    $ws.Columns.Item($c).ColumnWidth = 7.15
}
